# Applies the cryptocurrency table refresh described in the commit
# "Updated cryptos list on Sat Nov 11 15:47:56 UTC 2023 with GitHub Actions".
# All target cells are plain text (prices / percentages rendered as strings,
# some of which look numeric, e.g. "252.53"), so we force a text number format
# before assigning each value to avoid Excel auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.160.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.071.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.86"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +20.37%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.80"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.385"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +10.58%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.49"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.375.45"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +9.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.069.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.186.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0929"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +13.48%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +14.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.52"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.78"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.80%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.82"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.78"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.03%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.41"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.46%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.110"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +24.85%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.11"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.45"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +30.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.70"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +18.01%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +11.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.305.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.53%  "

Write-Host "Applied 100 cell updates"
